$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$data = @(
    ,@(2, 'Xiaomi Redmi Watch 5 Active, matte silver - Smartwatch', 39.99, 'Euronics')
    ,@(3, 'Xiaomi Redmi Buds 5, black - Wireless earbuds', 39.99, 'Euronics')
    ,@(4, 'Xiaomi Redmi Watch 5 Lite, black - Smartwatch', 49.99, 'Euronics')
    ,@(5, 'Xiaomi Redmi Watch 5 Lite, light gold - Smartwatch', 49.99, 'Euronics')
    ,@(6, 'Redmi A5 3/ 64GB Lake Green Xiaomi', 76.98999999999999, 'RD Electronics')
    ,@(7, 'Redmi A5 3/ 64GB Midnight Black Xiaomi', 79, 'RD Electronics')
    ,@(8, 'Redmi A5 3/ 64GB Sandy Gold Xiaomi', 79, 'RD Electronics')
    ,@(9, 'Redmi A5 3/ 64GB Ocean Blue Xiaomi', 79, 'RD Electronics')
    ,@(10, 'Redmi A5 4/ 128GB Ocean Blue Xiaomi', 94.98999999999999, 'RD Electronics')
    ,@(11, 'Redmi A5 4/ 128GB Lake Green Xiaomi', 94.98999999999999, 'RD Electronics')
    ,@(12, 'Redmi A5 4/ 128GB Midnight Black Xiaomi', 94.98999999999999, 'RD Electronics')
    ,@(13, 'Xiaomi Redmi Watch 5, purple - Smartwatch', 109.99, 'Euronics')
    ,@(14, 'Xiaomi Redmi Watch 5, black - Smartwatch', 109.99, 'Euronics')
    ,@(15, 'Xiaomi Redmi A5, midnight black - Smartphone', 109.99, 'Euronics')
    ,@(16, 'Xiaomi Redmi A5, sandy gold - Smartphone', 109.99, 'Euronics')
    ,@(17, 'Xiaomi Redmi Watch 5, silver - Smartwatch', 109.99, 'Euronics')
    ,@(18, 'Xiaomi Redmi A5, ocean blue - Smartphone', 109.99, 'Euronics')
    ,@(19, 'Redmi A5 4/ 128GB Sandy Gold Xiaomi', 124.99, 'RD Electronics')
    ,@(20, 'Xiaomi Redmi Note 12 5G, 128 GB, green - Smartphone', 179.99, 'Euronics')
    ,@(21, 'Xiaomi Redmi Note 13 5G, 256 GB, blue - Smartphone', 219.99, 'Euronics')
    ,@(22, 'Xiaomi Redmi Note 14 5G, 6 GB, 128 GB, midnight black - Smartphone', 229.99, 'Euronics')
    ,@(23, 'Xiaomi Redmi Note 14 5G, 6 GB, 128 GB, coral green - Smartphone', 229.99, 'Euronics')
    ,@(24, 'Xiaomi Redmi Note 14 5G, 8 GB, 256 GB, coral green - Smartphone', 249.99, 'Euronics')
    ,@(25, 'Xiaomi Redmi Note 14 5G, 8 GB, 256 GB, midnight black - Smartphone', 249.99, 'Euronics')
    ,@(26, 'Xiaomi Redmi Note 14 Pro 5G, lavender purple - Smartphone', 299.99, 'Euronics')
    ,@(27, 'Xiaomi Redmi Note 14 Pro 5G, coral green - Smartphone', 299.99, 'Euronics')
    ,@(28, 'Xiaomi Redmi Note 14 Pro 5G, midnight black - Smartphone', 299.99, 'Euronics')
    ,@(29, 'Xiaomi Redmi Note 14 Pro+ 5G, frost blue - Smartphone', 399.99, 'Euronics')
    ,@(30, 'Xiaomi Redmi Note 14 Pro+ 5G, midnight black - Smartphone', 399.99, 'Euronics')
    ,@(31, 'Xiaomi Redmi Note 14 Pro+ 5G, lavender purple - Smartphone', 399.99, 'Euronics')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
